$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.350872874259949
$ws.Range("B1").Value = 2.87433934211731
$ws.Range("C1").Value = 3.983532190322876
$ws.Range("D1").Value = 3.239475250244141
$ws.Range("E1").Value = 0.7739195227622986
